$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumo_por_Cliente")

$ws.Cells.Item(2, 10).Value = "INATIVO - 55.6 meses sem comprar"
$ws.Cells.Item(4, 10).Value = "INATIVO - 36.3 meses sem comprar"
$ws.Cells.Item(5, 10).Value = "INATIVO - 15.1 meses sem comprar"
$ws.Cells.Item(6, 10).Value = "INATIVO - 16.6 meses sem comprar"
$ws.Cells.Item(8, 10).Value = "INATIVO - 18.3 meses sem comprar"
$ws.Cells.Item(9, 10).Value = "INATIVO - 19.3 meses sem comprar"
$ws.Cells.Item(10, 10).Value = "INATIVO - 2.8 meses sem comprar"
$ws.Cells.Item(11, 10).Value = "INATIVO - 5.2 meses sem comprar"
$ws.Cells.Item(12, 10).Value = "INATIVO - 0.3 meses sem comprar"
$ws.Cells.Item(16, 10).Value = "INATIVO - 40.1 meses sem comprar"
$ws.Cells.Item(17, 10).Value = "INATIVO - 2.7 meses sem comprar"
$ws.Cells.Item(18, 10).Value = "INATIVO - 37.7 meses sem comprar"
$ws.Cells.Item(19, 10).Value = "INATIVO - 11.6 meses sem comprar"
$ws.Cells.Item(20, 10).Value = "INATIVO - 15.6 meses sem comprar"
$ws.Cells.Item(21, 10).Value = "INATIVO - 37.5 meses sem comprar"
$ws.Cells.Item(23, 10).Value = "INATIVO - 29.4 meses sem comprar"
$ws.Cells.Item(24, 10).Value = "INATIVO - 38.1 meses sem comprar"
$ws.Cells.Item(25, 10).Value = "INATIVO - 0.1 meses sem comprar"
$ws.Cells.Item(27, 10).Value = "INATIVO - 19.8 meses sem comprar"
$ws.Cells.Item(29, 10).Value = "INATIVO - 17.9 meses sem comprar"
$ws.Cells.Item(30, 10).Value = "INATIVO - 6.3 meses sem comprar"
$ws.Cells.Item(31, 10).Value = "INATIVO - 7.1 meses sem comprar"
$ws.Cells.Item(32, 10).Value = "INATIVO - 22.8 meses sem comprar"
$ws.Cells.Item(33, 10).Value = "INATIVO - 7.8 meses sem comprar"
$ws.Cells.Item(35, 10).Value = "INATIVO - 14.9 meses sem comprar"
$ws.Cells.Item(36, 10).Value = "INATIVO - 27.0 meses sem comprar"
$ws.Cells.Item(37, 10).Value = "INATIVO - 6.7 meses sem comprar"
$ws.Cells.Item(39, 10).Value = "INATIVO - 32.5 meses sem comprar"
$ws.Cells.Item(40, 10).Value = "INATIVO - 34.6 meses sem comprar"
$ws.Cells.Item(41, 10).Value = "INATIVO - 12.8 meses sem comprar"
$ws.Cells.Item(43, 10).Value = "INATIVO - 7.7 meses sem comprar"
$ws.Cells.Item(44, 10).Value = "INATIVO - 26.0 meses sem comprar"
$ws.Cells.Item(45, 10).Value = "INATIVO - 6.4 meses sem comprar"
$ws.Cells.Item(46, 10).Value = "INATIVO - 16.2 meses sem comprar"
$ws.Cells.Item(47, 10).Value = "INATIVO - 2.8 meses sem comprar"
$ws.Cells.Item(49, 10).Value = "INATIVO - 7.0 meses sem comprar"
$ws.Cells.Item(50, 10).Value = "INATIVO - 7.8 meses sem comprar"
$ws.Cells.Item(51, 10).Value = "INATIVO - 10.2 meses sem comprar"
$ws.Cells.Item(52, 10).Value = "INATIVO - 6.4 meses sem comprar"
$ws.Cells.Item(57, 5).Value = 32
$ws.Cells.Item(57, 8).Value = 45842.80699074074
$ws.Cells.Item(57, 9).Value = 45857.80699074074
$ws.Cells.Item(58, 5).Value = 11
$ws.Cells.Item(58, 8).Value = 45844.76849537037
$ws.Cells.Item(58, 9).Value = 45906.76849537037
$ws.Cells.Item(59, 10).Value = "INATIVO - 12.2 meses sem comprar"
$ws.Cells.Item(65, 10).Value = "INATIVO - 28.1 meses sem comprar"
$ws.Cells.Item(66, 10).Value = "INATIVO - 22.1 meses sem comprar"
$ws.Cells.Item(67, 2).Value = 0.33
$ws.Cells.Item(67, 3).Value = 0.33
$ws.Cells.Item(67, 5).Value = 32
$ws.Cells.Item(67, 8).Value = 45843.95056712963
$ws.Cells.Item(67, 9).Value = 45874.95056712963
$ws.Cells.Item(68, 10).Value = "INATIVO - 12.8 meses sem comprar"
$ws.Cells.Item(69, 10).Value = "INATIVO - 11.6 meses sem comprar"
$ws.Cells.Item(71, 3).Value = 0.17
$ws.Cells.Item(71, 4).Value = 0.67
$ws.Cells.Item(71, 5).Value = 14
$ws.Cells.Item(71, 6).Value = 0.67
$ws.Cells.Item(71, 8).Value = 45843.9522337963
$ws.Cells.Item(71, 9).Value = 45874.9522337963
$ws.Cells.Item(72, 2).Value = 0.25
$ws.Cells.Item(72, 3).Value = 0.17
$ws.Cells.Item(72, 4).Value = 0.33
$ws.Cells.Item(72, 5).Value = 12
$ws.Cells.Item(72, 6).Value = 0.33
$ws.Cells.Item(72, 7).Value = "1x por mês - irregular (preferencialmente na 1ª quinzena)"
$ws.Cells.Item(72, 8).Value = 45844.89819444445
$ws.Cells.Item(72, 9).Value = 45875.89819444445
$ws.Cells.Item(74, 10).Value = "INATIVO - 22.1 meses sem comprar"
$ws.Cells.Item(75, 10).Value = "INATIVO - 33.7 meses sem comprar"
$ws.Cells.Item(76, 10).Value = "INATIVO - 8.0 meses sem comprar"
$ws.Cells.Item(79, 10).Value = "INATIVO - 6.4 meses sem comprar"
$ws.Cells.Item(82, 10).Value = "INATIVO - 22.5 meses sem comprar"
$ws.Cells.Item(84, 10).Value = "INATIVO - 26.5 meses sem comprar"
$ws.Cells.Item(85, 10).Value = "INATIVO - 6.2 meses sem comprar"
$ws.Cells.Item(86, 10).Value = "INATIVO - 22.0 meses sem comprar"
$ws.Cells.Item(87, 10).Value = "INATIVO - 9.7 meses sem comprar"
$ws.Cells.Item(88, 10).Value = "INATIVO - 15.3 meses sem comprar"
$ws.Cells.Item(89, 10).Value = "INATIVO - 5.2 meses sem comprar"
$ws.Cells.Item(90, 10).Value = "INATIVO - 12.0 meses sem comprar"
$ws.Cells.Item(91, 10).Value = "INATIVO - 11.4 meses sem comprar"
$ws.Cells.Item(92, 10).Value = "INATIVO - 15.3 meses sem comprar"
$ws.Cells.Item(93, 10).Value = "INATIVO - 33.7 meses sem comprar"
$ws.Cells.Item(94, 10).Value = "INATIVO - 14.0 meses sem comprar"
$ws.Cells.Item(95, 10).Value = "INATIVO - 19.0 meses sem comprar"
$ws.Cells.Item(96, 10).Value = "INATIVO - 16.7 meses sem comprar"
$ws.Cells.Item(97, 10).Value = "INATIVO - 19.4 meses sem comprar"
$ws.Cells.Item(98, 10).Value = "INATIVO - 33.2 meses sem comprar"
$ws.Cells.Item(100, 10).Value = "INATIVO - 2.5 meses sem comprar"
$ws.Cells.Item(101, 10).Value = "INATIVO - 23.1 meses sem comprar"
$ws.Cells.Item(102, 10).Value = "INATIVO - 37.5 meses sem comprar"
$ws.Cells.Item(104, 10).Value = "INATIVO - 14.8 meses sem comprar"
$ws.Cells.Item(105, 10).Value = "INATIVO - 25.2 meses sem comprar"
$ws.Cells.Item(106, 10).Value = "INATIVO - 10.8 meses sem comprar"
$ws.Cells.Item(107, 10).Value = "INATIVO - 25.6 meses sem comprar"
$ws.Cells.Item(108, 10).Value = "INATIVO - 16.0 meses sem comprar"
$ws.Cells.Item(109, 10).Value = "INATIVO - 6.9 meses sem comprar"
$ws.Cells.Item(110, 10).Value = "INATIVO - 21.6 meses sem comprar"
$ws.Cells.Item(111, 10).Value = "INATIVO - 7.3 meses sem comprar"
$ws.Cells.Item(113, 10).Value = "INATIVO - 8.8 meses sem comprar"
$ws.Cells.Item(114, 5).Value = 16315
$ws.Cells.Item(114, 8).Value = 45842.74878472222
$ws.Cells.Item(114, 9).Value = 45843.74878472222
